$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.823.98"
$ws.Range("E2").Value = "  -0.12%  "
$ws.Range("D3").Value = "2.295.06"
$ws.Range("E3").Value = "  +0.02%  "
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "114.59"
$ws.Range("E5").Value = "  +17.60%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "269.79"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.626"
$ws.Range("E7").Value = "  +0.14%  "
$ws.Range("E8").Value = "  +0.25%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.617"
$ws.Range("E9").Value = "  +1.34%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "48.43"
$ws.Range("E10").Value = "  +6.68%  "
$ws.Range("E11").Value = "  +0.78%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.95"
$ws.Range("E12").Value = "  +13.55%  "
$ws.Range("E13").Value = "  +0.22%  "
$ws.Range("E14").Value = "  -0.22%  "
$ws.Range("D15").Value = "2.636.99"
$ws.Range("E15").Value = "  -0.04%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.857"
$ws.Range("E16").Value = "  -0.41%  "
$ws.Range("D17").Value = "2.291.35"
$ws.Range("E17").Value = "  -0.22%  "
$ws.Range("D18").Value = "43.715.21"
$ws.Range("E18").Value = "  -0.19%  "
$ws.Range("E19").Value = "  -0.78%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.99"
$ws.Range("E20").Value = "  +12.66%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.32"
$ws.Range("E21").Value = "  +0.10%  "
$ws.Range("E22").Value = "  -1.32%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.90"
$ws.Range("E23").Value = "  +8.89%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "232.86"
$ws.Range("E24").Value = "  -0.22%  "
$ws.Range("E25").Value = "  +7.95%  "
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.56"
$ws.Range("E27").Value = "  +2.16%  "
$ws.Range("B28").Value = "InjectiveProtocol"
$ws.Range("C28").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "42.59"
$ws.Range("E28").Value = "  +11.04%  "
$ws.Range("B29").Value = "LEO"
$ws.Range("C29").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.90"
$ws.Range("E29").Value = "  -0.32%  "
$ws.Range("E30").Value = "  -1.98%  "
$ws.Range("E31").Value = "  -1.40%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "175.48"
$ws.Range("E32").Value = "  -0.61%  "
$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0931"
$ws.Range("E33").Value = "  +4.30%  "
$ws.Range("B34").Value = "EthereumClassic"
$ws.Range("C34").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "21.62"
$ws.Range("E34").Value = "  -1.01%  "
$ws.Range("E35").Value = "  +5.35%  "
$ws.Range("E36").Value = "  -0.19%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.77"
$ws.Range("E37").Value = "  +1.18%  "
$ws.Range("E38").Value = "  +2.96%  "
$ws.Range("E39").Value = "  -3.39%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.82"
$ws.Range("E40").Value = "  +9.98%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "13.99"
$ws.Range("E41").Value = "  +14.71%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "74.67"
$ws.Range("E42").Value = "  +15.55%  "
$ws.Range("E43").Value = "  +2.70%  "
$ws.Range("E44").Value = "  +2.74%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "6.36"
$ws.Range("E45").Value = "  +21.72%  "
$ws.Range("E46").Value = "  +0.21%  "
$ws.Range("E47").Value = "  +2.32%  "
$ws.Range("E48").Value = "  +0.65%  "
$ws.Range("B49").Value = "TrustWalletToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.26"
$ws.Range("E49").Value = "  +2.40%  "
$ws.Range("B50").Value = "Aave"
$ws.Range("C50").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "102.39"
$ws.Range("E50").Value = "  +3.20%  "
$ws.Range("E51").Value = "  -2.96%  "
